# Update standard-score lookup tables (column B) across worksheets 1-20
# per the commit "TODC grade norms through iwr"
$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(4, 2).Value = 73
$ws.Cells.Item(5, 2).Value = 77
$ws.Cells.Item(6, 2).Value = 81
$ws.Cells.Item(7, 2).Value = 85
$ws.Cells.Item(8, 2).Value = 89
$ws.Cells.Item(9, 2).Value = 93
$ws.Cells.Item(10, 2).Value = 97
$ws.Cells.Item(11, 2).Value = 102
$ws.Cells.Item(12, 2).Value = 106
$ws.Cells.Item(14, 2).Value = 114
$ws.Cells.Item(15, 2).Value = 119
$ws.Cells.Item(16, 2).Value = 123
$ws.Cells.Item(17, 2).Value = 128

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 2).Value = 63
$ws.Cells.Item(3, 2).Value = 66
$ws.Cells.Item(4, 2).Value = 70
$ws.Cells.Item(7, 2).Value = 82
$ws.Cells.Item(8, 2).Value = 86
$ws.Cells.Item(12, 2).Value = 102
$ws.Cells.Item(13, 2).Value = 106
$ws.Cells.Item(16, 2).Value = 118
$ws.Cells.Item(17, 2).Value = 123
$ws.Cells.Item(18, 2).Value = 127

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 2).Value = 61
$ws.Cells.Item(3, 2).Value = 64
$ws.Cells.Item(4, 2).Value = 68
$ws.Cells.Item(5, 2).Value = 71
$ws.Cells.Item(12, 2).Value = 98
$ws.Cells.Item(13, 2).Value = 102
$ws.Cells.Item(17, 2).Value = 118
$ws.Cells.Item(18, 2).Value = 122
$ws.Cells.Item(19, 2).Value = 126

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 2).Value = 59
$ws.Cells.Item(3, 2).Value = 62
$ws.Cells.Item(4, 2).Value = 66
$ws.Cells.Item(5, 2).Value = 69
$ws.Cells.Item(6, 2).Value = 73
$ws.Cells.Item(7, 2).Value = 76
$ws.Cells.Item(9, 2).Value = 83
$ws.Cells.Item(16, 2).Value = 109
$ws.Cells.Item(17, 2).Value = 113
$ws.Cells.Item(18, 2).Value = 117
$ws.Cells.Item(19, 2).Value = 121
$ws.Cells.Item(20, 2).Value = 126

# Sheet 5
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(2, 2).Value = 57
$ws.Cells.Item(3, 2).Value = 61
$ws.Cells.Item(4, 2).Value = 64
$ws.Cells.Item(5, 2).Value = 67
$ws.Cells.Item(6, 2).Value = 71
$ws.Cells.Item(7, 2).Value = 74
$ws.Cells.Item(8, 2).Value = 77
$ws.Cells.Item(10, 2).Value = 84
$ws.Cells.Item(12, 2).Value = 91
$ws.Cells.Item(14, 2).Value = 98
$ws.Cells.Item(15, 2).Value = 102
$ws.Cells.Item(18, 2).Value = 113
$ws.Cells.Item(19, 2).Value = 117
$ws.Cells.Item(20, 2).Value = 121
$ws.Cells.Item(21, 2).Value = 125
$ws.Cells.Item(22, 2).Value = 129

# Sheet 6
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2, 2).Value = 56
$ws.Cells.Item(3, 2).Value = 59
$ws.Cells.Item(4, 2).Value = 62
$ws.Cells.Item(5, 2).Value = 65
$ws.Cells.Item(6, 2).Value = 69
$ws.Cells.Item(7, 2).Value = 72
$ws.Cells.Item(8, 2).Value = 75
$ws.Cells.Item(9, 2).Value = 78
$ws.Cells.Item(10, 2).Value = 82
$ws.Cells.Item(11, 2).Value = 85
$ws.Cells.Item(14, 2).Value = 95
$ws.Cells.Item(16, 2).Value = 102
$ws.Cells.Item(17, 2).Value = 106
$ws.Cells.Item(18, 2).Value = 110
$ws.Cells.Item(19, 2).Value = 113
$ws.Cells.Item(20, 2).Value = 117
$ws.Cells.Item(21, 2).Value = 121
$ws.Cells.Item(22, 2).Value = 125
$ws.Cells.Item(23, 2).Value = 129

# Sheet 7
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(2, 2).Value = 55
$ws.Cells.Item(3, 2).Value = 58
$ws.Cells.Item(4, 2).Value = 61
$ws.Cells.Item(5, 2).Value = 64
$ws.Cells.Item(6, 2).Value = 67
$ws.Cells.Item(7, 2).Value = 70
$ws.Cells.Item(8, 2).Value = 73
$ws.Cells.Item(9, 2).Value = 76
$ws.Cells.Item(10, 2).Value = 79
$ws.Cells.Item(11, 2).Value = 82
$ws.Cells.Item(12, 2).Value = 85
$ws.Cells.Item(13, 2).Value = 89
$ws.Cells.Item(14, 2).Value = 92
$ws.Cells.Item(15, 2).Value = 95
$ws.Cells.Item(16, 2).Value = 99
$ws.Cells.Item(17, 2).Value = 102
$ws.Cells.Item(18, 2).Value = 106
$ws.Cells.Item(19, 2).Value = 109
$ws.Cells.Item(20, 2).Value = 113
$ws.Cells.Item(21, 2).Value = 116
$ws.Cells.Item(22, 2).Value = 120
$ws.Cells.Item(23, 2).Value = 124
$ws.Cells.Item(24, 2).Value = 128

# Sheet 8
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(2, 2).Value = 53
$ws.Cells.Item(3, 2).Value = 56
$ws.Cells.Item(4, 2).Value = 59
$ws.Cells.Item(5, 2).Value = 62
$ws.Cells.Item(6, 2).Value = 64
$ws.Cells.Item(7, 2).Value = 67
$ws.Cells.Item(8, 2).Value = 70
$ws.Cells.Item(9, 2).Value = 73
$ws.Cells.Item(10, 2).Value = 76
$ws.Cells.Item(11, 2).Value = 79
$ws.Cells.Item(12, 2).Value = 82
$ws.Cells.Item(13, 2).Value = 85
$ws.Cells.Item(14, 2).Value = 89
$ws.Cells.Item(15, 2).Value = 92
$ws.Cells.Item(16, 2).Value = 95
$ws.Cells.Item(17, 2).Value = 98
$ws.Cells.Item(18, 2).Value = 101
$ws.Cells.Item(19, 2).Value = 105
$ws.Cells.Item(20, 2).Value = 108
$ws.Cells.Item(21, 2).Value = 112
$ws.Cells.Item(22, 2).Value = 115
$ws.Cells.Item(23, 2).Value = 119
$ws.Cells.Item(24, 2).Value = 123
$ws.Cells.Item(25, 2).Value = 126

# Sheet 9
$ws = $wb.Worksheets.Item(9)
$ws.Cells.Item(2, 2).Value = 52
$ws.Cells.Item(3, 2).Value = 55
$ws.Cells.Item(4, 2).Value = 57
$ws.Cells.Item(5, 2).Value = 60
$ws.Cells.Item(6, 2).Value = 63
$ws.Cells.Item(7, 2).Value = 65
$ws.Cells.Item(8, 2).Value = 68
$ws.Cells.Item(9, 2).Value = 71
$ws.Cells.Item(10, 2).Value = 74
$ws.Cells.Item(11, 2).Value = 77
$ws.Cells.Item(12, 2).Value = 80
$ws.Cells.Item(13, 2).Value = 83
$ws.Cells.Item(14, 2).Value = 85
$ws.Cells.Item(15, 2).Value = 88
$ws.Cells.Item(16, 2).Value = 92
$ws.Cells.Item(17, 2).Value = 95
$ws.Cells.Item(18, 2).Value = 98
$ws.Cells.Item(19, 2).Value = 101
$ws.Cells.Item(20, 2).Value = 104
$ws.Cells.Item(21, 2).Value = 107
$ws.Cells.Item(22, 2).Value = 111
$ws.Cells.Item(23, 2).Value = 114
$ws.Cells.Item(24, 2).Value = 118
$ws.Cells.Item(25, 2).Value = 121
$ws.Cells.Item(26, 2).Value = 125
$ws.Cells.Item(27, 2).Value = 129

# Sheet 10
$ws = $wb.Worksheets.Item(10)
$ws.Cells.Item(2, 2).Value = 51
$ws.Cells.Item(3, 2).Value = 53
$ws.Cells.Item(4, 2).Value = 56
$ws.Cells.Item(5, 2).Value = 59
$ws.Cells.Item(6, 2).Value = 61
$ws.Cells.Item(7, 2).Value = 64
$ws.Cells.Item(8, 2).Value = 66
$ws.Cells.Item(9, 2).Value = 69
$ws.Cells.Item(10, 2).Value = 72
$ws.Cells.Item(11, 2).Value = 74
$ws.Cells.Item(12, 2).Value = 77
$ws.Cells.Item(13, 2).Value = 80
$ws.Cells.Item(14, 2).Value = 83
$ws.Cells.Item(15, 2).Value = 86
$ws.Cells.Item(16, 2).Value = 89
$ws.Cells.Item(17, 2).Value = 92
$ws.Cells.Item(18, 2).Value = 95
$ws.Cells.Item(19, 2).Value = 98
$ws.Cells.Item(20, 2).Value = 101
$ws.Cells.Item(21, 2).Value = 104
$ws.Cells.Item(22, 2).Value = 107
$ws.Cells.Item(23, 2).Value = 110
$ws.Cells.Item(24, 2).Value = 114
$ws.Cells.Item(25, 2).Value = 117
$ws.Cells.Item(26, 2).Value = 121
$ws.Cells.Item(27, 2).Value = 124
$ws.Cells.Item(28, 2).Value = 128

# Sheet 11
$ws = $wb.Worksheets.Item(11)
$ws.Cells.Item(2, 2).Value = 50
$ws.Cells.Item(3, 2).Value = 52
$ws.Cells.Item(4, 2).Value = 55
$ws.Cells.Item(5, 2).Value = 57
$ws.Cells.Item(6, 2).Value = 60
$ws.Cells.Item(7, 2).Value = 62
$ws.Cells.Item(8, 2).Value = 65
$ws.Cells.Item(9, 2).Value = 67
$ws.Cells.Item(10, 2).Value = 70
$ws.Cells.Item(11, 2).Value = 73
$ws.Cells.Item(12, 2).Value = 75
$ws.Cells.Item(13, 2).Value = 78
$ws.Cells.Item(14, 2).Value = 81
$ws.Cells.Item(15, 2).Value = 83
$ws.Cells.Item(16, 2).Value = 86
$ws.Cells.Item(17, 2).Value = 89
$ws.Cells.Item(18, 2).Value = 92
$ws.Cells.Item(19, 2).Value = 95
$ws.Cells.Item(20, 2).Value = 98
$ws.Cells.Item(21, 2).Value = 101
$ws.Cells.Item(22, 2).Value = 104
$ws.Cells.Item(23, 2).Value = 107
$ws.Cells.Item(24, 2).Value = 110
$ws.Cells.Item(25, 2).Value = 113
$ws.Cells.Item(26, 2).Value = 117
$ws.Cells.Item(27, 2).Value = 120
$ws.Cells.Item(28, 2).Value = 124
$ws.Cells.Item(29, 2).Value = 127

# Sheet 12
$ws = $wb.Worksheets.Item(12)
$ws.Cells.Item(2, 2).Value = 49
$ws.Cells.Item(3, 2).Value = 52
$ws.Cells.Item(4, 2).Value = 54
$ws.Cells.Item(5, 2).Value = 56
$ws.Cells.Item(6, 2).Value = 59
$ws.Cells.Item(7, 2).Value = 61
$ws.Cells.Item(8, 2).Value = 64
$ws.Cells.Item(9, 2).Value = 66
$ws.Cells.Item(10, 2).Value = 68
$ws.Cells.Item(11, 2).Value = 71
$ws.Cells.Item(12, 2).Value = 73
$ws.Cells.Item(13, 2).Value = 76
$ws.Cells.Item(14, 2).Value = 79
$ws.Cells.Item(15, 2).Value = 81
$ws.Cells.Item(16, 2).Value = 84
$ws.Cells.Item(17, 2).Value = 87
$ws.Cells.Item(18, 2).Value = 89
$ws.Cells.Item(19, 2).Value = 92
$ws.Cells.Item(20, 2).Value = 95
$ws.Cells.Item(21, 2).Value = 98
$ws.Cells.Item(22, 2).Value = 101
$ws.Cells.Item(23, 2).Value = 104
$ws.Cells.Item(24, 2).Value = 107
$ws.Cells.Item(25, 2).Value = 110
$ws.Cells.Item(26, 2).Value = 113
$ws.Cells.Item(27, 2).Value = 116
$ws.Cells.Item(28, 2).Value = 120
$ws.Cells.Item(29, 2).Value = 123
$ws.Cells.Item(30, 2).Value = 127

# Sheet 13
$ws = $wb.Worksheets.Item(13)
$ws.Cells.Item(2, 2).Value = 49
$ws.Cells.Item(3, 2).Value = 51
$ws.Cells.Item(4, 2).Value = 53
$ws.Cells.Item(5, 2).Value = 56
$ws.Cells.Item(6, 2).Value = 58
$ws.Cells.Item(7, 2).Value = 60
$ws.Cells.Item(8, 2).Value = 62
$ws.Cells.Item(9, 2).Value = 65
$ws.Cells.Item(10, 2).Value = 67
$ws.Cells.Item(11, 2).Value = 70
$ws.Cells.Item(12, 2).Value = 72
$ws.Cells.Item(13, 2).Value = 74
$ws.Cells.Item(14, 2).Value = 77
$ws.Cells.Item(15, 2).Value = 79
$ws.Cells.Item(16, 2).Value = 82
$ws.Cells.Item(17, 2).Value = 85
$ws.Cells.Item(18, 2).Value = 87
$ws.Cells.Item(19, 2).Value = 90
$ws.Cells.Item(20, 2).Value = 92
$ws.Cells.Item(21, 2).Value = 95
$ws.Cells.Item(22, 2).Value = 98
$ws.Cells.Item(23, 2).Value = 101
$ws.Cells.Item(24, 2).Value = 104
$ws.Cells.Item(25, 2).Value = 107
$ws.Cells.Item(26, 2).Value = 110
$ws.Cells.Item(27, 2).Value = 113
$ws.Cells.Item(28, 2).Value = 117
$ws.Cells.Item(29, 2).Value = 120
$ws.Cells.Item(30, 2).Value = 124
$ws.Cells.Item(31, 2).Value = 127

# Sheet 14
$ws = $wb.Worksheets.Item(14)
$ws.Cells.Item(2, 2).Value = 48
$ws.Cells.Item(3, 2).Value = 51
$ws.Cells.Item(4, 2).Value = 53
$ws.Cells.Item(5, 2).Value = 55
$ws.Cells.Item(6, 2).Value = 57
$ws.Cells.Item(7, 2).Value = 59
$ws.Cells.Item(8, 2).Value = 62
$ws.Cells.Item(9, 2).Value = 64
$ws.Cells.Item(10, 2).Value = 66
$ws.Cells.Item(11, 2).Value = 68
$ws.Cells.Item(12, 2).Value = 71
$ws.Cells.Item(13, 2).Value = 73
$ws.Cells.Item(14, 2).Value = 75
$ws.Cells.Item(15, 2).Value = 78
$ws.Cells.Item(16, 2).Value = 80
$ws.Cells.Item(17, 2).Value = 83
$ws.Cells.Item(18, 2).Value = 85
$ws.Cells.Item(19, 2).Value = 88
$ws.Cells.Item(20, 2).Value = 90
$ws.Cells.Item(21, 2).Value = 93
$ws.Cells.Item(22, 2).Value = 96
$ws.Cells.Item(23, 2).Value = 99
$ws.Cells.Item(24, 2).Value = 101
$ws.Cells.Item(25, 2).Value = 104
$ws.Cells.Item(26, 2).Value = 107
$ws.Cells.Item(27, 2).Value = 110
$ws.Cells.Item(28, 2).Value = 114
$ws.Cells.Item(29, 2).Value = 117
$ws.Cells.Item(30, 2).Value = 120
$ws.Cells.Item(31, 2).Value = 124
$ws.Cells.Item(32, 2).Value = 128

# Sheet 15
$ws = $wb.Worksheets.Item(15)
$ws.Cells.Item(2, 2).Value = 48
$ws.Cells.Item(3, 2).Value = 50
$ws.Cells.Item(4, 2).Value = 52
$ws.Cells.Item(5, 2).Value = 54
$ws.Cells.Item(6, 2).Value = 56
$ws.Cells.Item(7, 2).Value = 59
$ws.Cells.Item(8, 2).Value = 61
$ws.Cells.Item(9, 2).Value = 63
$ws.Cells.Item(10, 2).Value = 65
$ws.Cells.Item(11, 2).Value = 67
$ws.Cells.Item(12, 2).Value = 70
$ws.Cells.Item(13, 2).Value = 72
$ws.Cells.Item(14, 2).Value = 74
$ws.Cells.Item(15, 2).Value = 76
$ws.Cells.Item(16, 2).Value = 79
$ws.Cells.Item(17, 2).Value = 81
$ws.Cells.Item(18, 2).Value = 84
$ws.Cells.Item(19, 2).Value = 86
$ws.Cells.Item(20, 2).Value = 89
$ws.Cells.Item(21, 2).Value = 91
$ws.Cells.Item(22, 2).Value = 94
$ws.Cells.Item(23, 2).Value = 96
$ws.Cells.Item(24, 2).Value = 99
$ws.Cells.Item(25, 2).Value = 102
$ws.Cells.Item(26, 2).Value = 105
$ws.Cells.Item(27, 2).Value = 108
$ws.Cells.Item(28, 2).Value = 111
$ws.Cells.Item(29, 2).Value = 114
$ws.Cells.Item(30, 2).Value = 117
$ws.Cells.Item(31, 2).Value = 121
$ws.Cells.Item(32, 2).Value = 125
$ws.Cells.Item(33, 2).Value = 128

# Sheet 16
$ws = $wb.Worksheets.Item(16)
$ws.Cells.Item(2, 2).Value = 48
$ws.Cells.Item(3, 2).Value = 50
$ws.Cells.Item(4, 2).Value = 52
$ws.Cells.Item(5, 2).Value = 54
$ws.Cells.Item(6, 2).Value = 56
$ws.Cells.Item(7, 2).Value = 58
$ws.Cells.Item(8, 2).Value = 60
$ws.Cells.Item(9, 2).Value = 62
$ws.Cells.Item(10, 2).Value = 64
$ws.Cells.Item(11, 2).Value = 66
$ws.Cells.Item(12, 2).Value = 69
$ws.Cells.Item(13, 2).Value = 71
$ws.Cells.Item(14, 2).Value = 73
$ws.Cells.Item(15, 2).Value = 75
$ws.Cells.Item(16, 2).Value = 78
$ws.Cells.Item(17, 2).Value = 80
$ws.Cells.Item(18, 2).Value = 82
$ws.Cells.Item(20, 2).Value = 87
$ws.Cells.Item(21, 2).Value = 89
$ws.Cells.Item(22, 2).Value = 92
$ws.Cells.Item(24, 2).Value = 97
$ws.Cells.Item(25, 2).Value = 100
$ws.Cells.Item(26, 2).Value = 103
$ws.Cells.Item(29, 2).Value = 112
$ws.Cells.Item(30, 2).Value = 115
$ws.Cells.Item(31, 2).Value = 118
$ws.Cells.Item(32, 2).Value = 122
$ws.Cells.Item(33, 2).Value = 126

# Sheet 17
$ws = $wb.Worksheets.Item(17)
$ws.Cells.Item(2, 2).Value = 48
$ws.Cells.Item(3, 2).Value = 50
$ws.Cells.Item(4, 2).Value = 52
$ws.Cells.Item(5, 2).Value = 53
$ws.Cells.Item(6, 2).Value = 55
$ws.Cells.Item(7, 2).Value = 57
$ws.Cells.Item(8, 2).Value = 59
$ws.Cells.Item(9, 2).Value = 61
$ws.Cells.Item(10, 2).Value = 63
$ws.Cells.Item(11, 2).Value = 65
$ws.Cells.Item(12, 2).Value = 67
$ws.Cells.Item(13, 2).Value = 70
$ws.Cells.Item(14, 2).Value = 72
$ws.Cells.Item(15, 2).Value = 74
$ws.Cells.Item(16, 2).Value = 76
$ws.Cells.Item(17, 2).Value = 78
$ws.Cells.Item(18, 2).Value = 80
$ws.Cells.Item(20, 2).Value = 85
$ws.Cells.Item(21, 2).Value = 87
$ws.Cells.Item(23, 2).Value = 92
$ws.Cells.Item(25, 2).Value = 97
$ws.Cells.Item(26, 2).Value = 100
$ws.Cells.Item(29, 2).Value = 108
$ws.Cells.Item(32, 2).Value = 118

# Sheet 18
$ws = $wb.Worksheets.Item(18)
$ws.Cells.Item(2, 2).Value = 48
$ws.Cells.Item(3, 2).Value = 50
$ws.Cells.Item(4, 2).Value = 51
$ws.Cells.Item(5, 2).Value = 53
$ws.Cells.Item(6, 2).Value = 55
$ws.Cells.Item(7, 2).Value = 57
$ws.Cells.Item(8, 2).Value = 59
$ws.Cells.Item(9, 2).Value = 61
$ws.Cells.Item(10, 2).Value = 62
$ws.Cells.Item(11, 2).Value = 64
$ws.Cells.Item(12, 2).Value = 66
$ws.Cells.Item(13, 2).Value = 68
$ws.Cells.Item(14, 2).Value = 70
$ws.Cells.Item(15, 2).Value = 72
$ws.Cells.Item(16, 2).Value = 74
$ws.Cells.Item(17, 2).Value = 76
$ws.Cells.Item(18, 2).Value = 78
$ws.Cells.Item(22, 2).Value = 87
$ws.Cells.Item(25, 2).Value = 94
$ws.Cells.Item(28, 2).Value = 102
$ws.Cells.Item(32, 2).Value = 114
$ws.Cells.Item(34, 2).Value = 121
$ws.Cells.Item(35, 2).Value = 125

# Sheet 19
$ws = $wb.Worksheets.Item(19)
$ws.Cells.Item(2, 2).Value = 48
$ws.Cells.Item(11, 2).Value = 64
$ws.Cells.Item(15, 2).Value = 71
$ws.Cells.Item(16, 2).Value = 73
$ws.Cells.Item(17, 2).Value = 75
$ws.Cells.Item(18, 2).Value = 77
$ws.Cells.Item(21, 2).Value = 83
$ws.Cells.Item(22, 2).Value = 85
$ws.Cells.Item(26, 2).Value = 94
$ws.Cells.Item(31, 2).Value = 106
$ws.Cells.Item(32, 2).Value = 109
$ws.Cells.Item(34, 2).Value = 116
$ws.Cells.Item(36, 2).Value = 124

# Sheet 20
$ws = $wb.Worksheets.Item(20)
$ws.Cells.Item(2, 2).Value = 50
$ws.Cells.Item(3, 2).Value = 51
$ws.Cells.Item(4, 2).Value = 53
$ws.Cells.Item(5, 2).Value = 54
$ws.Cells.Item(6, 2).Value = 56
$ws.Cells.Item(7, 2).Value = 57
$ws.Cells.Item(8, 2).Value = 59
$ws.Cells.Item(9, 2).Value = 60
$ws.Cells.Item(10, 2).Value = 62
$ws.Cells.Item(11, 2).Value = 63
$ws.Cells.Item(12, 2).Value = 65
$ws.Cells.Item(13, 2).Value = 67
$ws.Cells.Item(14, 2).Value = 68
$ws.Cells.Item(15, 2).Value = 70
$ws.Cells.Item(16, 2).Value = 72
$ws.Cells.Item(17, 2).Value = 73
$ws.Cells.Item(18, 2).Value = 75
$ws.Cells.Item(19, 2).Value = 77
$ws.Cells.Item(20, 2).Value = 79
$ws.Cells.Item(21, 2).Value = 81
$ws.Cells.Item(22, 2).Value = 83
$ws.Cells.Item(23, 2).Value = 84
$ws.Cells.Item(26, 2).Value = 91
$ws.Cells.Item(27, 2).Value = 93
$ws.Cells.Item(31, 2).Value = 102
$ws.Cells.Item(32, 2).Value = 105
$ws.Cells.Item(33, 2).Value = 108
$ws.Cells.Item(34, 2).Value = 111
$ws.Cells.Item(35, 2).Value = 114
$ws.Cells.Item(36, 2).Value = 118
$ws.Cells.Item(37, 2).Value = 122
$ws.Cells.Item(38, 2).Value = 127

